$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "População"
$ws.Range("B3").Value = "PIB"
$ws.Range("B4").Value = "PIB per capta"
$ws.Range("B5").Value = "Altitude"
$ws.Range("B6").Value = "Acesso à Tratamento de Água"
$ws.Range("B7").Value = "Acesso à Coleta de Lixo"
$ws.Range("B8").Value = "Porcetagem de População Urbana"
